# Quiz 5 answers, and the Mid term review
$d = $word.ActiveDocument

# --- Q: "When is the objects constructor called?" -> answer (para 4 is
#     currently an empty paragraph right after the question)
$p = $d.Paragraphs(4)
$p.Range.InsertAfter("When the object tis instantiated/created. ")
$p.Range.Font.Color = 255

# --- Q: "When is the objects destructor called?" -> answer (para 6 is
#     currently an empty ListParagraph-styled paragraph; drop the list
#     style back to Normal so only the red run color remains)
$p = $d.Paragraphs(6)
$p.Style = "Normal"
$p.Range.InsertAfter("When the object is destroyed, by either reaching the end of its scope, or when it is deallocated.")
$p.Range.Font.Color = 255

# --- Split "What concept of object orientated programming..." right after
#     "orien" and drop the editing-cursor (_GoBack) bookmark there. Word
#     only keeps one _GoBack bookmark, so this also removes it from its old
#     spot at the very end of the document.
$p = $d.Paragraphs(8)
$splitPos = $p.Range.Start + "What concept of object orien".Length
$cursor = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $cursor) | Out-Null

# --- Q: "...make data members private..." -> answer "Encapsulation"
#     (para 9 is currently empty)
$p = $d.Paragraphs(9)
$p.Range.InsertAfter("Encapsulation")
$p.Range.Font.Color = 255

# --- Q: "What state should an object be in after the constructor..." ->
#     answer (para 12 is currently an empty ListParagraph paragraph with a
#     360-twip left indent; drop both back to Normal so only the red run
#     color remains)
$p = $d.Paragraphs(12)
$p.Style = "Normal"
$p.Range.InsertAfter("A safe state.")
$p.Range.Font.Color = 255
